# Slide 6 ("Reference - ..."), Content Placeholder 2: the third paragraph
# ("Hey look, I was browsing through the list of libraries ... cryptography.")
# originally had its opening sentence split across three separate <a:r> runs:
#   "Hey look, I " + "was browsing " + "through the list of libraries on the
#   left and found a built-in library just for "
# Collapse those three runs back into a single run (keeping the formatting
# of the first one), matching the upstream "fixed hyperlinks" commit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$run1Text = "Hey look, I "
$run2Text = "was browsing "
$run3Text = "through the list of libraries on the left and found a built-in library just for "

$full = $tr.Text
$startIdx = $full.IndexOf($run1Text + $run2Text + $run3Text)
if ($startIdx -lt 0) {
    throw "Could not locate the target run sequence in the shape's text."
}

# TextRange.Characters(Start, Length) uses 1-based character offsets.
$run1Start = $startIdx + 1
$run2Start = $run1Start + $run1Text.Length
$run3Start = $run2Start + $run2Text.Length

# Edit from the last run back to the first so earlier offsets stay valid
# while later ones are being collapsed to empty strings.
$c3 = $tr.Characters($run3Start, $run3Text.Length)
$c3.Text = ""

$c2 = $tr.Characters($run2Start, $run2Text.Length)
$c2.Text = ""

$c1 = $tr.Characters($run1Start, $run1Text.Length)
$c1.Text = $run1Text + $run2Text + $run3Text
